# feat: add 2022-Q1 data
#
# - the existing "总计" sheet is turned into the new "2022-Q1" quarterly
#   fund-holding report (same column layout as "2021-Q3" / "2021-Q4")
# - a fresh "总计" sheet is appended with the running totals table, now
#   including the new 2022-Q1 row on top

$wb = $excel.ActiveWorkbook

$oldTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1) Duplicate the current "总计" sheet first (it already carries the
#    right sheetPr/format/pageMargins for the totals table) and place the
#    duplicate at the very end - that duplicate will become the refreshed
#    "总计" sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$oldTotal.Copy($null, $lastSheet)
$newTotal = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------------
# 2) Turn the original "总计" sheet into "2022-Q1": same report layout as
#    the other quarterly sheets (基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#    持有市值(亿元)/仓位排名). Rename it out of the way first so the
#    duplicate can reclaim the "总计" name.
# ---------------------------------------------------------------------
$q1 = $oldTotal
$q1.Name = "2022-Q1"
$newTotal.Name = "总计"

# drop the old 3rd totals row - the report sheet only has one data row
$q1.Rows.Item(3).Delete()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"

# extend the header row with the extra report columns, re-using D1's
# existing header style for the new cells
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# data row - 基金规模/股票总仓位/仓位占比/持有市值(亿元) are stored as plain text
# (matching every other quarter sheet), 仓位排名 is numeric
$q1Text = $q1.Range("B2:G2")
$q1Text.NumberFormat = "@"
$q1.Range("B2").Value = "968013"
$q1.Range("C2").Value = "施罗德亚洲高息股债基金M"
$q1.Range("D2").Value = "297.64"
$q1.Range("E2").Value = "57.54"
$q1.Range("F2").Value = "1.66"
$q1.Range("G2").Value = "4.9408"
$q1Text.Style = "Normal"
$q1.Range("H2").Value = 3

# ---------------------------------------------------------------------
# 3) Refresh the new "总计" sheet: add the 2022-Q1 row on top, pushing the
#    2021-Q4 / 2021-Q3 rows down by one.
# ---------------------------------------------------------------------
$newTotal.Range("A2").Copy()
$newTotal.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 1
$newTotal.Range("D2").Value = 4.94

$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 1
$newTotal.Range("D3").Value = 5.03

$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q3"
$newTotal.Range("C4").Value = 1
$newTotal.Range("D4").Value = 8.34

# ---------------------------------------------------------------------
# 4) Restore "2021-Q3" as the selected tab (unchanged by this edit).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
